$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that originally sits
#    right after the "Play All Lucky Clover Slot for Free" heading.
$metaOld = 'Meta description: Read our review of All Lucky Clover online slot game and play for free with expanding Wilds and Scatters paying out regardless of position.'
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $metaOld) {
        $null = $p.Range.Delete()
        break
    }
}

# 2. Insert a new bold "Play All Lucky Clover Slot for Free" paragraph
#    right before the trailing "Feature image prompt for DALLE" paragraph.
$n = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($n)
$null = $last.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($n)
$frag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play All Lucky Clover Slot for Free</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $newPara.Range.InsertXML($frag)

# 3. Swap out the old DALLE feature-image prompt text for the meta
#    description text, keeping the existing (italic) run formatting.
$old = 'Feature image prompt for DALLE: Create a colorful cartoon-style image for "All Lucky Clover" slot game that reflects the game''s fun and upbeat theme. The image should feature a happy Maya warrior wearing glasses. The warrior can be surrounded by a field of clovers or holding a clover in their hand. The illustration should be bright and cheerful, with a mix of greens, golds, and other bold colors that complement the game''s overall aesthetics. The title of the game should be prominently displayed in the image, along with some of the game''s symbols, such as fruits, the clover jewel, horseshoe, and diamond. The image should be eye-catching, inviting, and representative of the game''s exciting features and potential payouts.'
$new = 'Read our review of All Lucky Clover online slot game and play for free with expanding Wilds and Scatters paying out regardless of position.'
$null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
